$wb = $excel.ActiveWorkbook

# --- Demand sheet: extend the single future-year demand value into a full
#     12-step time series (t=1..12), replacing the old single high value
#     with the new annualised value repeated across the horizon. ---
$wsDemand = $wb.Worksheets.Item("Demand")

# Row 3 (t=1) keeps its format but gets the new value.
$wsDemand.Range("B3").Value = 531666667

# Add rows 4-14 (t=2..12) mirroring row 3: same value, same number format.
for ($t = 2; $t -le 12; $t++) {
    $row = $t + 2
    $wsDemand.Range("A" + $row).Value = $t
    $wsDemand.Range("B" + $row).Value = 531666667
}

# Copy row 3's formatting down across the newly-added B cells so they pick
# up the same (#,##0-style) number format as the existing data.
$wsDemand.Range("B3").Copy()
$wsDemand.Range("B4:B14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Widen column B so the larger numbers aren't truncated.
$wsDemand.Columns.Item(2).ColumnWidth = 10.166666666666666

# --- Selection / active-sheet bookkeeping: the workbook was left with the
#     SupIm sheet's selection parked at A15, and "Demand" as the active tab
#     with H14 selected. ---
$wsSupIm = $wb.Worksheets.Item("SupIm")
$wsSupIm.Activate() | Out-Null
$wsSupIm.Range("A15").Select() | Out-Null

$wsDemand.Activate() | Out-Null
$wsDemand.Range("H14").Select() | Out-Null
